$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 668.6429000000001
$ws.Range("I53").Value = 840.4
$ws.Range("J53").Value = 573.2222
$ws.Range("K53").Value = 840.4
$ws.Range("L53").Value = 573.2222
$ws.Range("M53").Value = -203.4
$ws.Range("N53").Value = -1847.2222
$ws.Range("H58").Value = 215
$ws.Range("I58").Value = 215
$ws.Range("K58").Value = 645
$ws.Range("M58").Value = -495
$ws.Range("H112").Value = 1288.375
$ws.Range("I112").Value = 684.2
$ws.Range("J112").Value = 1447.3684
$ws.Range("K112").Value = 2052.6
$ws.Range("L112").Value = 4342.1052
$ws.Range("M112").Value = -944.6000000000004
$ws.Range("N112").Value = -6558.1052
$ws.Range("H132").Value = 1360.4082
$ws.Range("I132").Value = 1003.3333
$ws.Range("K132").Value = 3009.9999
$ws.Range("M132").Value = -479.9998999999998
$ws.Range("H134").Value = 54992.285
$ws.Range("J134").Value = 54992.285
$ws.Range("L134").Value = 54992.285
$ws.Range("N134").Value = -65132.285

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 18002560
$ws.Range("I45").Value = 3729.75
$ws.Range("K45").Value = 3729.75
$ws.Range("M45").Value = -3352.75
$ws.Range("H74").Value = 58021.61
$ws.Range("I74").Value = 85336.836
$ws.Range("J74").Value = 3391.1667
$ws.Range("K74").Value = 85336.836
$ws.Range("L74").Value = 3391.1667
$ws.Range("M74").Value = -84462.836
$ws.Range("N74").Value = -5139.1667
$ws.Range("H77").Value = 58021.61
$ws.Range("I77").Value = 85336.836
$ws.Range("J77").Value = 3391.1667
$ws.Range("K77").Value = 426684.18
$ws.Range("L77").Value = 16955.8335
$ws.Range("M77").Value = -422316.18
$ws.Range("N77").Value = -25691.8335
$ws.Range("H122").Value = 1871.6923
$ws.Range("I122").Value = 1788.3
$ws.Range("J122").Value = 2149.6667
$ws.Range("K122").Value = 5364.9
$ws.Range("L122").Value = 6449.000100000001
$ws.Range("M122").Value = -2914.9
$ws.Range("N122").Value = -11349.0001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3298.0356
$ws.Range("J86").Value = 5595.5713
$ws.Range("L86").Value = 5595.5713
$ws.Range("N86").Value = -7841.5713
$ws.Range("H89").Value = 3298.0356
$ws.Range("J89").Value = 5595.5713
$ws.Range("L89").Value = 27977.8565
$ws.Range("N89").Value = -39209.85649999999
$ws.Range("M133").Value = -19940
$ws.Range("H133").Value = 25000
$ws.Range("I133").Value = 25000
$ws.Range("K133").Value = 25000

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2441.4666
$ws.Range("I31").Value = 1872.2609
$ws.Range("K31").Value = 1872.2609
$ws.Range("M31").Value = -1577.2609
$ws.Range("H34").Value = 2441.4666
$ws.Range("I34").Value = 1872.2609
$ws.Range("K34").Value = 1872.2609
$ws.Range("M34").Value = -1670.2609
$ws.Range("H132").Value = 2068949.4
$ws.Range("I132").Value = 2843681
$ws.Range("K132").Value = 8531043
$ws.Range("M132").Value = -8528513
$ws.Range("H134").Value = 2589931.8
$ws.Range("I134").Value = 3762974.2
$ws.Range("J134").Value = 113509.22
$ws.Range("K134").Value = 11288922.6
$ws.Range("L134").Value = 340527.66
$ws.Range("M134").Value = -11286387.6
$ws.Range("N134").Value = -345597.66

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 169832.67
$ws.Range("I14").Value = 169832.67
$ws.Range("K14").Value = 509498.01
$ws.Range("M14").Value = -509325.01
$ws.Range("N58").ClearContents()
$ws.Range("H58").Value = 2299.5
$ws.Range("I58").Value = 2299.5
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 6898.5
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -6770.5
$ws.Range("H80").Value = 2000
$ws.Range("J80").Value = 2000
$ws.Range("L80").Value = 6000
$ws.Range("N80").Value = -7872
$ws.Range("H83").Value = 2000
$ws.Range("J83").Value = 2000
$ws.Range("L83").Value = 18000
$ws.Range("N83").Value = -27360
$ws.Range("H113").Value = 3242199.2
$ws.Range("I113").Value = 1337.25
$ws.Range("J113").Value = 4420694.5
$ws.Range("K113").Value = 4011.75
$ws.Range("L113").Value = 13262083.5
$ws.Range("M113").Value = -1841.75
$ws.Range("N113").Value = -13266423.5
$ws.Range("H121").Value = 2427.182
$ws.Range("I121").Value = 599.8570999999999
$ws.Range("J121").Value = 5625
$ws.Range("K121").Value = 1799.5713
$ws.Range("L121").Value = 16875
$ws.Range("M121").Value = -489.5712999999998
$ws.Range("N121").Value = -19495

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1144.7778
$ws.Range("I107").Value = 1061.1666
$ws.Range("K107").Value = 1061.1666
$ws.Range("M107").Value = 858.8334
$ws.Range("N122").Value = -16900
$ws.Range("H122").Value = 915179.5600000001
$ws.Range("I122").Value = 1117663.9
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 3352991.7
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -3350541.7

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2538.4443
$ws.Range("I82").Value = 2174.3333
$ws.Range("J82").Value = 3266.6667
$ws.Range("K82").Value = 2174.3333
$ws.Range("L82").Value = 3266.6667
$ws.Range("M82").Value = -1813.3333
$ws.Range("N82").Value = -3988.6667
$ws.Range("H85").Value = 2538.4443
$ws.Range("I85").Value = 2174.3333
$ws.Range("J85").Value = 3266.6667
$ws.Range("K85").Value = 2174.3333
$ws.Range("L85").Value = 3266.6667
$ws.Range("M85").Value = -926.3332999999998
$ws.Range("N85").Value = -5762.6667
$ws.Range("N130").ClearContents()
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("H132").Value = 1709.2354
$ws.Range("I132").Value = 1659.8125
$ws.Range("K132").Value = 4979.4375
$ws.Range("M132").Value = -2449.4375

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 525001.5
$ws.Range("I2").Value = 50000
$ws.Range("K2").Value = 50000
$ws.Range("M2").Value = -49888
$ws.Range("H126").Value = 2944.08
$ws.Range("I126").Value = 2397
$ws.Range("J126").Value = 3916.6667
$ws.Range("K126").Value = 7191
$ws.Range("L126").Value = 11750.0001
$ws.Range("M126").Value = -4721
$ws.Range("N126").Value = -16690.0001
$ws.Range("H136").Value = 1353
$ws.Range("I136").Value = 1141
$ws.Range("J136").Value = 2201
$ws.Range("K136").Value = 3423
$ws.Range("L136").Value = 6603
$ws.Range("M136").Value = -873
$ws.Range("N136").Value = -11703
